# Finished 10mm flexor testing: updated the last calibration reading
# (10 mm Festo flexor) from 4.89 V to 4.82 V, and left the sheet's
# active-cell selection on C8 (as last left by the author).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the final voltage reading in the calibration table.
$ws.Range("A6").Value = 4.82

# Move the active selection to C8.
$ws.Range("C8").Select() | Out-Null
